{"js": "// Update the three-digit x one-digit multiplication answers in the table.\n// Each old value is unique in the document, so a direct search+replace\n// (matching the whole cell text) is unambiguous.\nconst replacements = [\n  [\"341\u00d79=3069\", \"320\u00d79=2880\"],\n  [\"659\u00d78=5272\", \"660\u00d75=3300\"],\n  [\"387\u00d76=2322\", \"833\u00d79=7497\"],\n  [\"399\u00d76=2394\", \"817\u00d75=4085\"],\n  [\"185\u00d79=1665\", \"351\u00d72=702\"],\n  [\"746\u00d74=2984\", \"937\u00d79=8433\"],\n  [\"670\u00d79=6030\", \"217\u00d76=1302\"],\n  [\"860\u00d74=3440\", \"667\u00d79=6003\"],\n  [\"422\u00d73=1266\", \"427\u00d73=1281\"],\n  [\"672\u00d79=6048\", \"983\u00d75=4915\"],\n  [\"851\u00d77=5957\", \"298\u00d77=2086\"],\n  [\"966\u00d78=7728\", \"795\u00d74=3180\"],\n  [\"401\u00d72=802\", \"986\u00d73=2958\"],\n  [\"721\u00d74=2884\", \"484\u00d72=968\"],\n  [\"765\u00d77=5355\", \"217\u00d72=434\"],\n  [\"500\u00d78=4000\", \"186\u00d72=372\"],\n  [\"321\u00d76=1926\", \"495\u00d76=2970\"],\n  [\"147\u00d77=1029\", \"770\u00d74=3080\"],\n  [\"140\u00d73=420\", \"559\u00d72=1118\"],\n  [\"284\u00d75=1420\", \"696\u00d75=3480\"],\n  [\"664\u00d74=2656\", \"621\u00d77=4347\"],\n  [\"292\u00d77=2044\", \"350\u00d79=3150\"],\n  [\"366\u00d74=1464\", \"994\u00d79=8946\"],\n  [\"141\u00d77=987\", \"377\u00d74=1508\"],\n  [\"358\u00d78=2864\", \"455\u00d78=3640\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the three-digit x one-digit multiplication answers in the table.\n# Each old value is unique in the document, so a direct Find/Replace\n# (matching the whole cell text) is unambiguous.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"341\u00d79=3069\", \"320\u00d79=2880\"),\n    @(\"659\u00d78=5272\", \"660\u00d75=3300\"),\n    @(\"387\u00d76=2322\", \"833\u00d79=7497\"),\n    @(\"399\u00d76=2394\", \"817\u00d75=4085\"),\n    @(\"185\u00d79=1665\", \"351\u00d72=702\"),\n    @(\"746\u00d74=2984\", \"937\u00d79=8433\"),\n    @(\"670\u00d79=6030\", \"217\u00d76=1302\"),\n    @(\"860\u00d74=3440\", \"667\u00d79=6003\"),\n    @(\"422\u00d73=1266\", \"427\u00d73=1281\"),\n    @(\"672\u00d79=6048\", \"983\u00d75=4915\"),\n    @(\"851\u00d77=5957\", \"298\u00d77=2086\"),\n    @(\"966\u00d78=7728\", \"795\u00d74=3180\"),\n    @(\"401\u00d72=802\", \"986\u00d73=2958\"),\n    @(\"721\u00d74=2884\", \"484\u00d72=968\"),\n    @(\"765\u00d77=5355\", \"217\u00d72=434\"),\n    @(\"500\u00d78=4000\", \"186\u00d72=372\"),\n    @(\"321\u00d76=1926\", \"495\u00d76=2970\"),\n    @(\"147\u00d77=1029\", \"770\u00d74=3080\"),\n    @(\"140\u00d73=420\", \"559\u00d72=1118\"),\n    @(\"284\u00d75=1420\", \"696\u00d75=3480\"),\n    @(\"664\u00d74=2656\", \"621\u00d77=4347\"),\n    @(\"292\u00d77=2044\", \"350\u00d79=3150\"),\n    @(\"366\u00d74=1464\", \"994\u00d79=8946\"),\n    @(\"141\u00d77=987\", \"377\u00d74=1508\"),\n    @(\"358\u00d78=2864\", \"455\u00d78=3640\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
